# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" data snapshot: new timestamp, updated case
# counts for several countries, and two countries (Tunez, Guatemala) whose
# rank (by total cases) moved up enough to shift several neighbouring rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refreshed "as of" timestamp
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 19:22"

# Row 4: Estados Unidos - updated totals
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 300148
$ws.Range("C4").Value = 22987
$ws.Range("D4").Value = 14464
$ws.Range("E4").Value = 277543
$ws.Range("F4").Value = 6333
$ws.Range("G4").Value = 737
$ws.Range("H4").Value = 8141

# Row 16: Canada - updated totals
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 12956
$ws.Range("C16").Value = 581
$ws.Range("D16").Value = 2322
$ws.Range("E16").Value = 10420
$ws.Range("F16").Value = 120
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 214

# Row 37: Pakistan - updated totals
$ws.Range("A37").Value = "Pakistan"
$ws.Range("B37").Value = 2818
$ws.Range("C37").Value = 132
$ws.Range("D37").Value = 131
$ws.Range("E37").Value = 2646
$ws.Range("F37").Value = 13
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 41

# Rows 73-79: Tunez jumps ahead of Kazajistan/Azerbaiyan/Libano/Letonia/
# Camerun/Bulgaria, pushing each of them down one row
$ws.Range("A73").Value = "Tunez"
$ws.Range("B73").Value = 553
$ws.Range("C73").Value = 58
$ws.Range("D73").Value = 5
$ws.Range("E73").Value = 529
$ws.Range("F73").Value = 26
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 19

$ws.Range("A74").Value = "Kazajistan"
$ws.Range("B74").Value = 525
$ws.Range("C74").Value = 61
$ws.Range("D74").Value = 36
$ws.Range("E74").Value = 484
$ws.Range("F74").Value = 6
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 5

$ws.Range("A75").Value = "Azerbaiyan"
$ws.Range("B75").Value = 521
$ws.Range("C75").Value = 78
$ws.Range("D75").Value = 32
$ws.Range("E75").Value = 484
$ws.Range("F75").Value = 17
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 5

$ws.Range("A76").Value = "Libano"
$ws.Range("B76").Value = 520
$ws.Range("C76").Value = 12
$ws.Range("D76").Value = 54
$ws.Range("E76").Value = 449
$ws.Range("F76").Value = 26
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 17

$ws.Range("A77").Value = "Letonia"
$ws.Range("B77").Value = 509
$ws.Range("C77").Value = 16
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 507
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 1

$ws.Range("A78").Value = "Camerun"
$ws.Range("B78").Value = 509
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 17
$ws.Range("E78").Value = 484
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 8

$ws.Range("A79").Value = "Bulgaria"
$ws.Range("B79").Value = 503
$ws.Range("C79").Value = 18
$ws.Range("D79").Value = 34
$ws.Range("E79").Value = 452
$ws.Range("F79").Value = 26
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 17

# Row 85: Republica de Chipre - updated totals
$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("B85").Value = 426
$ws.Range("C85").Value = 30
$ws.Range("D85").Value = 33
$ws.Range("E85").Value = 382
$ws.Range("F85").Value = 11
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 11

# Row 109: Sri Lanka - updated totals
$ws.Range("A109").Value = "Sri Lanka"
$ws.Range("B109").Value = 166
$ws.Range("C109").Value = 7
$ws.Range("D109").Value = 27
$ws.Range("E109").Value = 134
$ws.Range("F109").Value = 5
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 5

# Rows 133-138: Guatemala jumps ahead of Guayana Francesa/El Salvador/
# Jamaica/Barbados/Republica de Yibuti, pushing each of them down one row
$ws.Range("A133").Value = "Guatemala"
$ws.Range("B133").Value = 57
$ws.Range("C133").Value = 7
$ws.Range("D133").Value = 15
$ws.Range("E133").Value = 40
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 2

$ws.Range("A134").Value = "Guayana Francesa"
$ws.Range("B134").Value = 57
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 22
$ws.Range("E134").Value = 35
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

$ws.Range("A135").Value = "El Salvador"
$ws.Range("B135").Value = 56
$ws.Range("C135").Value = 10
$ws.Range("D135").Value = 0
$ws.Range("E135").Value = 53
$ws.Range("F135").Value = 4
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 3

$ws.Range("A136").Value = "Jamaica"
$ws.Range("B136").Value = 53
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 7
$ws.Range("E136").Value = 43
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 3

$ws.Range("A137").Value = "Barbados"
$ws.Range("B137").Value = 51
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 51
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

$ws.Range("A138").Value = "Republica de Yibuti"
$ws.Range("B138").Value = 50
$ws.Range("C138").Value = 1
$ws.Range("D138").Value = 8
$ws.Range("E138").Value = 42
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0
